$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 16.363986
$ws.Range("D3").Value = 69.225762
